# Apply updated dSF (column F) values per repulled data / mean calculation fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 0
$ws.Range("F12").Value = 7
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 12
$ws.Range("F17").Value = 1
$ws.Range("F20").Value = -9
$ws.Range("F24").Value = -4
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 0
